$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# sheet1
$ws1.Range("F6").Value = 386
$ws1.Range("F7").Value = 1218
$ws1.Range("F8").Value = 457
$ws1.Range("F9").Value = 7424
$ws1.Range("F10").Value = 87
$ws1.Range("F11").Value = 98
$ws1.Range("F12").Value = 2062
$ws1.Range("F13").Value = 8078
$ws1.Range("F14").Value = 48
$ws1.Range("F16").Value = 5545
$ws1.Range("F18").Value = 2478
$ws1.Range("F19").Value = 1051
$ws1.Range("F21").Value = 309
$ws1.Range("F23").Value = 83
$ws1.Range("F25").Value = 412
$ws1.Range("F26").Value = 725
$ws1.Range("F27").Value = 17
$ws1.Range("F28").Value = 2516
$ws1.Range("F30").Value = 282
$ws1.Range("F31").Value = 94
$ws1.Range("F32").Value = 192
$ws1.Range("F33").Value = 609
$ws1.Range("F36").Value = 1559
$ws1.Range("F39").Value = 2455
$ws1.Range("F40").Value = 2236
$ws1.Range("F41").Value = 4

# sheet2
$ws2.Range("F3").Value = 89
$ws2.Range("F4").Value = 84
$ws2.Range("F5").Value = 24
$ws2.Range("F6").Value = 30
$ws2.Range("F7").Value = 30

# sheet3
$ws3.Range("F3").Value = 1284

# sheet4
$ws4.Range("F4").Value = 1284
$ws4.Range("F6").Value = 386
$ws4.Range("F7").Value = 1218
$ws4.Range("F8").Value = 457
$ws4.Range("F9").Value = 7424
$ws4.Range("F10").Value = 87
$ws4.Range("F11").Value = 98
$ws4.Range("F12").Value = 2062
$ws4.Range("F13").Value = 8078
$ws4.Range("F14").Value = 48
$ws4.Range("F16").Value = 5545
$ws4.Range("F18").Value = 2478
$ws4.Range("F19").Value = 1051
$ws4.Range("F21").Value = 309
$ws4.Range("F23").Value = 83
$ws4.Range("F24").Value = 89
$ws4.Range("F26").Value = 84
$ws4.Range("F27").Value = 412
$ws4.Range("F28").Value = 725
$ws4.Range("F29").Value = 17
$ws4.Range("F30").Value = 2516
$ws4.Range("F32").Value = 282
$ws4.Range("F33").Value = 94
$ws4.Range("F34").Value = 192
$ws4.Range("F35").Value = 24
$ws4.Range("F36").Value = 609
$ws4.Range("F39").Value = 30
$ws4.Range("F40").Value = 1559
$ws4.Range("F43").Value = 2455
$ws4.Range("F44").Value = 30
$ws4.Range("F45").Value = 2236
$ws4.Range("F46").Value = 4
